# testdata.xlsx update:
#  - "random" sheet accumulates newly generated sample e-mail addresses
#    (it already held "emailAddress"/"foxh@test.com" in A1:A2); append the
#    next batch of generated addresses below the existing data.
#  - "Registration" sheet's sample row (A2) is refreshed to use the newest
#    generated address from that batch.

$wb = $excel.ActiveWorkbook

$randomSheet = $wb.Worksheets.Item("random")
$randomSheet.Range("A3").Value = "ypqh@test.com"
$randomSheet.Range("A4").Value = "tdia@test.com"
$randomSheet.Range("A5").Value = "kxbh@test.com"
$randomSheet.Range("A6").Value = "ohge@test.com"

$registrationSheet = $wb.Worksheets.Item("Registration")
$registrationSheet.Range("A2").Value = "wjhk@test.com"
